$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- Row 15 ---
# Date column (A) must stay plain text (like all other date cells in this
# sheet), not get auto-converted to a serial date number. Pre-formatting
# as Text forces that, then we restore the default "Standard" cell style
# afterwards so the cell ends up unformatted, just like its neighbours.
$ws.Range("A15:A16").NumberFormat = "@"

$ws.Range("A15").Value = "01.01.2024"
$ws.Range("B15").Value = 135
$ws.Range("C15").Value = "Further worked on request handling"
$ws.Range("D15").Value = "21:00-23:15"

# --- Row 16 ---
$ws.Range("A16").Value = "02.01.2023"
$ws.Range("C16").Value = "Further worked on request handling"
$ws.Range("D16").Value = "12:40-zeit"

# Restore column A cells to the workbook's default style (no special
# number format), matching every other date cell in the table.
$ws.Range("A15:A16").Style = "Standard"

# Column D keeps the same "time span" cell style used by the rows above it
# (copy formatting only, so the text we just entered is left untouched).
$ws.Range("D14").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("D16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Reflect the final selection left by the editor.
$ws.Range("D16").Select()

$wb.Save()
